$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" '64.607.63'
$ws.Range("E2").Value = '  +1.40%  '
Set-TextValue $ws "D3" '3.434.96'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.06%  '
Set-TextValue $ws "D5" '574.70'
$ws.Range("E5").Value = '  -0.58%  '
Set-TextValue $ws "D6" '159.78'
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("E7").Value = '  +0.05%  '
Set-TextValue $ws "D8" '3.439.88'
$ws.Range("E8").Value = '  -0.01%  '
Set-TextValue $ws "D9" '0.581'
$ws.Range("E9").Value = '  +8.64%  '
Set-TextValue $ws "D10" '7.35'
$ws.Range("E10").Value = '  -3.04%  '
Set-TextValue $ws "D11" '0.125'
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("E12").Value = '  +0.48%  '
Set-TextValue $ws "D13" '4.037.17'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  -2.40%  '
Set-TextValue $ws "D15" '0.0000193'
$ws.Range("E15").Value = '  +3.56%  '
Set-TextValue $ws "D16" '28.18'
$ws.Range("E16").Value = '  +3.19%  '
Set-TextValue $ws "D17" '64.639.36'
$ws.Range("E17").Value = '  +1.33%  '
Set-TextValue $ws "D18" '3.442.07'
$ws.Range("E18").Value = '  +0.25%  '
Set-TextValue $ws "D19" '6.34'
$ws.Range("E19").Value = '  -1.56%  '
Set-TextValue $ws "D20" '14.24'
$ws.Range("E20").Value = '  -0.56%  '
Set-TextValue $ws "D21" '385.42'
$ws.Range("E21").Value = '  -1.50%  '
Set-TextValue $ws "D22" '8.15'
$ws.Range("E22").Value = '  -4.06%  '
Set-TextValue $ws "D23" '73.17'
$ws.Range("E23").Value = '  +1.50%  '
Set-TextValue $ws "D24" '0.543'
$ws.Range("E24").Value = '  +0.46%  '
Set-TextValue $ws "D25" '0.999'
$ws.Range("E25").Value = '  +0.35%  '
Set-TextValue $ws "D26" '0.0000123'
$ws.Range("E26").Value = '  +13.24%  '
Set-TextValue $ws "D27" '9.74'
$ws.Range("E27").Value = '  +2.27%  '
Set-TextValue $ws "D28" '0.179'
$ws.Range("E28").Value = '  -1.18%  '
Set-TextValue $ws "D29" '0.998'
$ws.Range("E29").Value = '  -0.20%  '
Set-TextValue $ws "D30" '6.19'
$ws.Range("E30").Value = '  +6.91%  '
$ws.Range("E31").Value = '  +3.37%  '
Set-TextValue $ws "D32" '2.04'
$ws.Range("E32").Value = '  -0.33%  '
Set-TextValue $ws "D33" '23.64'
$ws.Range("E33").Value = '  +0.48%  '
Set-TextValue $ws "D34" '6.53'
$ws.Range("E34").Value = '  -3.10%  '
Set-TextValue $ws "D36" '7.07'
$ws.Range("E36").Value = '  +3.25%  '
Set-TextValue $ws "D37" '163.14'
$ws.Range("E37").Value = '  +2.79%  '
Set-TextValue $ws "D38" '1.50'
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("E39").Value = '  +1.27%  '
Set-TextValue $ws "D40" '3.006.38'
$ws.Range("E40").Value = '  +4.22%  '
$ws.Range("E41").Value = '  -2.43%  '
Set-TextValue $ws "D42" '27.18'
$ws.Range("E42").Value = '  -3.50%  '
Set-TextValue $ws "D43" '4.57'
$ws.Range("E43").Value = '  +4.26%  '
Set-TextValue $ws "D44" '42.73'
$ws.Range("E44").Value = '  +1.85%  '
Set-TextValue $ws "D45" '0.0315'
$ws.Range("E45").Value = '  -2.02%  '
Set-TextValue $ws "D46" '0.772'
$ws.Range("E46").Value = '  +0.20%  '
Set-TextValue $ws "D47" '24.66'
$ws.Range("E47").Value = '  +8.95%  '
$ws.Range("E48").Value = '  -0.42%  '
Set-TextValue $ws "D49" '0.876'
$ws.Range("E49").Value = '  +5.69%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws "D50" '6.61'
$ws.Range("E50").Value = '  +3.32%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws "D51" '2.17'
$ws.Range("E51").Value = '  +2.45%  '
